$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update repaymentstrategy value (B17) from "RBI (India)" to the new scenario value
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move the active selection to B17, matching the recorded sheet view state
$ws.Range("B17").Select()
